$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 3, columns A..CY (1..103) with "b", except AP3="bb" and CI3="bbb"
for ($col = 1; $col -le 103; $col++) {
    $value = "b"
    if ($col -eq 42) { $value = "bb" }   # AP3
    if ($col -eq 87) { $value = "bbb" }  # CI3
    $ws.Cells.Item(3, $col).Value2 = $value
}

# Update the view: scroll so CO1 is the top-left cell and CZ3 is selected
$ws.Range("CZ3").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 93
$win.ScrollRow = 1
